{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// Implements the edit described by the diff:\n//  - Shortens the first paragraph (acting as the page's <meta description>)\n//    to be under 160 characters, splitting it into 3 runs as produced by Word.\n//  - Restores/normalizes the \"Budgetty\" bio paragraph and the collaborators'\n//    names with spell-check proofErr markers (Word re-flowed these runs\n//    after its proofing pass), with identical visible text.\n//  - Adds a grammar proofErr marker around \"move\" in the card-options\n//    paragraph, without changing the visible text.\n//  - Removes the stray <w:lastRenderedPageBreak/> marker.\n//  - Adds spell-check proofErr markers around each \"Budgetty\" occurrence in\n//    the closing image-caption paragraphs, splitting runs accordingly.\n\nfunction wrapBodyOoxml(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + bodyXml + '</w:body>' +\n    '</w:document>' +\n    '</pkg:xmlData>' +\n    '</pkg:part>' +\n    '</pkg:package>'\n  );\n}\n\n// Replaces the \"content\" (everything but the paragraph mark) of a paragraph\n// with a freshly authored <w:p>...</w:p> fragment, preserving exact run /\n// proofErr structure.\nfunction setParagraphOoxml(paragraph, innerXml) {\n  const range = paragraph.getRange(Word.RangeLocation.content);\n  range.insertOoxml(wrapBodyOoxml('<w:p>' + innerXml + '</w:p>'), Word.InsertLocation.replace);\n}\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// --- Paragraph 0: shortened \"meta description\" text, split into 3 runs ---\nsetParagraphOoxml(\n  paragraphs.items[0],\n  '<w:r><w:t>A</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> standalone digital wallet that encourages more mindful spending habits. I designed and animated the UI for the circular screen</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> in the middle.</w:t></w:r>'\n);\n\n// --- Paragraph 3: \"Budgetty is a standalone ... Sean Lee.\" (Collaboration) ---\nsetParagraphOoxml(\n  paragraphs.items[3],\n  '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Budgetty</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> is a standalone digital wallet that encourages more mindful spending habits in a world of abstract payment methods. I designed and animated the UI for the circular screen which users interact with through the surrounding </w:t></w:r>' +\n    '<w:r><w:t>dial</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">. This project was created in collaboration with Eduardo </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Zanforlin</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Mautner, Leo </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Baek</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Sanghyuk</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Seo</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>, and Sean Lee.</w:t></w:r>'\n);\n\n// --- Paragraph 7: \"The card options menu with a delete and move button. ...\" ---\nsetParagraphOoxml(\n  paragraphs.items[7],\n  '<w:r><w:t xml:space=\"preserve\">The card options menu with a delete and </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>move</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> button. The corresponding card is highlighted </w:t></w:r>' +\n    '<w:r><w:t>with</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> a white circle.</w:t></w:r>'\n);\n\n// --- Paragraph 15: remove the <w:lastRenderedPageBreak/> marker ---\nsetParagraphOoxml(\n  paragraphs.items[15],\n  '<w:r><w:t>By clicking on a specific transaction in the list, the user is shown the card details, the amount spent, the exact time of the transaction, and the merchant involved.</w:t></w:r>'\n);\n\n// --- Paragraph 17: \"Exploded view of Budgetty. From top to bottom: ...\" ---\nsetParagraphOoxml(\n  paragraphs.items[17],\n  '<w:r><w:t xml:space=\"preserve\">Exploded view of </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Budgetty</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">. From top to bottom: input </w:t></w:r>' +\n    '<w:r><w:t>dial</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">, circular </w:t></w:r>' +\n    '<w:r><w:t>screen, motherboard, motherboard screws, top cover, cover screws, bottom cover, and the wireless charger.</w:t></w:r>'\n);\n\n// --- Paragraph 18: \"Three Budgetty devices floating in mid-air, ...\" ---\nsetParagraphOoxml(\n  paragraphs.items[18],\n  '<w:r><w:t xml:space=\"preserve\">Three </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Budgetty</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> devices floating in mid-air, each displaying a different interface.</w:t></w:r>'\n);\n\n// --- Paragraph 19: \"Top-down view of Budgetty. ...\" ---\nsetParagraphOoxml(\n  paragraphs.items[19],\n  '<w:r><w:t xml:space=\"preserve\">Top-down view of </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Budgetty</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>. The circular screen at the center of its pentagonal shell displays the home page with the menu and remaining balance.</w:t></w:r>'\n);\n\n// --- Paragraph 20: \"Someone holding a 3D-printed resin model of Budgetty ...\" ---\nsetParagraphOoxml(\n  paragraphs.items[20],\n  '<w:r><w:t xml:space=\"preserve\">Someone holding a 3D-printed resin model of </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Budgetty</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> in their left hand while turning the input dial with their right hand.</w:t></w:r>'\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is the open document.\n#\n# Implements the edit described by the diff:\n#  - Shortens the first paragraph (acting as the page's <meta description>)\n#    to be under 160 characters, splitting it into 3 runs as produced by Word.\n#  - Restores/normalizes the \"Budgetty\" bio paragraph and the collaborators'\n#    names with spell-check proofErr markers (Word re-flowed these runs\n#    after its proofing pass), with identical visible text.\n#  - Adds a grammar proofErr marker around \"move\" in the card-options\n#    paragraph, without changing the visible text.\n#  - Removes the stray <w:lastRenderedPageBreak/> marker.\n#  - Adds spell-check proofErr markers around each \"Budgetty\" occurrence in\n#    the closing image-caption paragraphs, splitting runs accordingly.\n#\n# NOTE: calls are written as \"Set-ParagraphOoxml $var1 $var2\" (plain\n# variables only) rather than inlining \".Item(N)\" / parenthesized\n# concatenations directly as call arguments, to sidestep a parser quirk in\n# this host where a `.Method(...)` first argument immediately followed by a\n# parenthesized second argument loses that second argument.\n\nfunction Set-ParagraphOoxml($paragraph, [string]$innerXml) {\n    $xml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n        '<pkg:xmlData>' +\n        '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n        '<w:body><w:p>' + $innerXml + '</w:p></w:body>' +\n        '</w:document>' +\n        '</pkg:xmlData>' +\n        '</pkg:part>' +\n        '</pkg:package>'\n    $paragraph.Range.InsertXML($xml)\n}\n\nfunction Set-LastParagraphOoxml($doc, [string]$innerXml) {\n    # InsertXML-ing a whole <w:p> into the very last paragraph of the\n    # document's range leaves behind a duplicate trailing empty paragraph\n    # (the original final paragraph mark survives, pushed after the newly\n    # inserted one) because that final mark sits outside the story's\n    # addressable range. Undo that by merging the stray trailing empty\n    # paragraph back into the one we just wrote.\n    $paras = $doc.Paragraphs\n    $countBefore = $paras.Count\n    $target = $paras.Item($countBefore)\n    Set-ParagraphOoxml $target $innerXml\n\n    $parasAfter = $doc.Paragraphs\n    $countAfter = $parasAfter.Count\n    if ($countAfter -gt $countBefore) {\n        $trailing = $parasAfter.Item($countAfter)\n        $tr = $trailing.Range\n        $mergeRange = $doc.Range($tr.Start - 1, $tr.End)\n        $mergeRange.Delete()\n    }\n}\n\n$d = $word.ActiveDocument\n$paras = $d.Paragraphs\n\n# --- Paragraph 1: shortened \"meta description\" text, split into 3 runs ---\n$p1 = $paras.Item(1)\n$inner1 = '<w:r><w:t>A</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> standalone digital wallet that encourages more mindful spending habits. I designed and animated the UI for the circular screen</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> in the middle.</w:t></w:r>'\nSet-ParagraphOoxml $p1 $inner1\n\n# --- Paragraph 4: \"Budgetty is a standalone ... Sean Lee.\" (Collaboration) ---\n$p4 = $paras.Item(4)\n$inner4 = '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Budgetty</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> is a standalone digital wallet that encourages more mindful spending habits in a world of abstract payment methods. I designed and animated the UI for the circular screen which users interact with through the surrounding </w:t></w:r>' +\n    '<w:r><w:t>dial</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">. This project was created in collaboration with Eduardo </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Zanforlin</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> Mautner, Leo </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Baek</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">, </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Sanghyuk</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Seo</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>, and Sean Lee.</w:t></w:r>'\nSet-ParagraphOoxml $p4 $inner4\n\n# --- Paragraph 8: \"The card options menu with a delete and move button. ...\" ---\n$p8 = $paras.Item(8)\n$inner8 = '<w:r><w:t xml:space=\"preserve\">The card options menu with a delete and </w:t></w:r>' +\n    '<w:proofErr w:type=\"gramStart\"/>' +\n    '<w:r><w:t>move</w:t></w:r>' +\n    '<w:proofErr w:type=\"gramEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> button. The corresponding card is highlighted </w:t></w:r>' +\n    '<w:r><w:t>with</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\"> a white circle.</w:t></w:r>'\nSet-ParagraphOoxml $p8 $inner8\n\n# --- Paragraph 16: remove the <w:lastRenderedPageBreak/> marker ---\n$p16 = $paras.Item(16)\n$inner16 = '<w:r><w:t>By clicking on a specific transaction in the list, the user is shown the card details, the amount spent, the exact time of the transaction, and the merchant involved.</w:t></w:r>'\nSet-ParagraphOoxml $p16 $inner16\n\n# --- Paragraph 18: \"Exploded view of Budgetty. From top to bottom: ...\" ---\n$p18 = $paras.Item(18)\n$inner18 = '<w:r><w:t xml:space=\"preserve\">Exploded view of </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Budgetty</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\">. From top to bottom: input </w:t></w:r>' +\n    '<w:r><w:t>dial</w:t></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">, circular </w:t></w:r>' +\n    '<w:r><w:t>screen, motherboard, motherboard screws, top cover, cover screws, bottom cover, and the wireless charger.</w:t></w:r>'\nSet-ParagraphOoxml $p18 $inner18\n\n# --- Paragraph 19: \"Three Budgetty devices floating in mid-air, ...\" ---\n$p19 = $paras.Item(19)\n$inner19 = '<w:r><w:t xml:space=\"preserve\">Three </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Budgetty</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> devices floating in mid-air, each displaying a different interface.</w:t></w:r>'\nSet-ParagraphOoxml $p19 $inner19\n\n# --- Paragraph 20: \"Top-down view of Budgetty. ...\" ---\n$p20 = $paras.Item(20)\n$inner20 = '<w:r><w:t xml:space=\"preserve\">Top-down view of </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Budgetty</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>. The circular screen at the center of its pentagonal shell displays the home page with the menu and remaining balance.</w:t></w:r>'\nSet-ParagraphOoxml $p20 $inner20\n\n# --- Paragraph 21: \"Someone holding a 3D-printed resin model of Budgetty ...\" ---\n$p21 = $paras.Item(21)\n$inner21 = '<w:r><w:t xml:space=\"preserve\">Someone holding a 3D-printed resin model of </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>Budgetty</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t xml:space=\"preserve\"> in their left hand while turning the input dial with their right hand.</w:t></w:r>'\nSet-ParagraphOoxml $p21 $inner21\n"}
